$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "66.289.26"
$ws.Range("E2").Value = "  -1.03%  "

# Row 3
$ws.Range("D3").Value = "3.070.47"

# Row 4
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").Value = "'574.69"
$ws.Range("E5").Value = "  -0.82%  "

# Row 6
$ws.Range("E6").Value = "  -1.26%  "

# Row 7
$ws.Range("E7").Value = "  +0.05%  "

# Row 8
$ws.Range("D8").Value = "3.068.01"
$ws.Range("E8").Value = "  -1.58%  "

# Row 9
$ws.Range("E9").Value = "  -2.29%  "

# Row 10
$ws.Range("E10").Value = "  -1.85%  "

# Row 12
$ws.Range("E12").Value = "  -2.87%  "

# Row 13
$ws.Range("D13").Value = "'0.0000239"
$ws.Range("E13").Value = "  -3.86%  "

# Row 14
$ws.Range("D14").Value = "'35.63"
$ws.Range("E14").Value = "  -4.10%  "

# Row 16
$ws.Range("D16").Value = "3.581.06"
$ws.Range("E16").Value = "  -1.57%  "

# Row 17
$ws.Range("D17").Value = "66.207.14"
$ws.Range("E17").Value = "  -1.09%  "

# Row 18
$ws.Range("D18").Value = "'6.94"
$ws.Range("E18").Value = "  -3.10%  "

# Row 19
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "3.071.99"
$ws.Range("E19").Value = "  -1.57%  "

# Row 20
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "'16.56"
$ws.Range("E20").Value = "  +1.60%  "

# Row 21
$ws.Range("D21").Value = "'484.81"
$ws.Range("E21").Value = "  +1.77%  "

# Row 22
$ws.Range("E22").Value = "  -3.52%  "

# Row 23
$ws.Range("E23").Value = "  -2.98%  "

# Row 24
$ws.Range("D24").Value = "'82.19"
$ws.Range("E24").Value = "  -1.98%  "

# Row 25
$ws.Range("D25").Value = "'12.61"
$ws.Range("E25").Value = "  -4.51%  "

# Row 26
$ws.Range("D26").Value = "'2.20"
$ws.Range("E26").Value = "  -3.49%  "

# Row 27
$ws.Range("E27").Value = "  -2.57%  "

# Row 28
$ws.Range("E28").Value = "  -0.15%  "

# Row 29
$ws.Range("D29").Value = "'7.80"
$ws.Range("E29").Value = "  -1.39%  "

# Row 30
$ws.Range("E30").Value = "  -5.53%  "

# Row 31
$ws.Range("E31").Value = "  -3.49%  "

# Row 32
$ws.Range("D32").Value = "'27.63"
$ws.Range("E32").Value = "  -3.38%  "

# Row 33
$ws.Range("E33").Value = "  -3.24%  "

# Row 34
$ws.Range("D34").Value = "0.0₃0915"
$ws.Range("E34").Value = "  -3.80%  "

# Row 35
$ws.Range("E35").Value = "  -0.06%  "

# Row 36
$ws.Range("D36").Value = "'47.87"
$ws.Range("E36").Value = "  +2.05%  "

# Row 37
$ws.Range("B37").Value = "Mantle"
$ws.Range("C37").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D37").Value = "'0.942"
$ws.Range("E37").Value = "  -3.46%  "

# Row 38
$ws.Range("B38").Value = "Filecoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D38").Value = "'5.55"
$ws.Range("E38").Value = "  -5.07%  "

# Row 39
$ws.Range("E39").Value = "  -1.16%  "

# Row 40
$ws.Range("B40").Value = "TheGraph"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D40").Value = "'0.300"
$ws.Range("E40").Value = "  -3.81%  "

# Row 41
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'1.95"
$ws.Range("E41").Value = "  -5.09%  "

# Row 42
$ws.Range("E42").Value = "  -4.57%  "

# Row 43
$ws.Range("D43").Value = "2.774.90"
$ws.Range("E43").Value = "  -1.78%  "

# Row 44
$ws.Range("D44").Value = "'2.53"
$ws.Range("E44").Value = "  -0.94%  "

# Row 45
$ws.Range("E45").Value = "  -2.83%  "

# Row 46
$ws.Range("D46").Value = "'134.46"
$ws.Range("E46").Value = "  -1.04%  "

# Row 47
$ws.Range("D47").Value = "'363.73"
$ws.Range("E47").Value = "  -4.96%  "

# Row 49
$ws.Range("D49").Value = "'24.18"
$ws.Range("E49").Value = "  -3.21%  "

# Row 50
$ws.Range("E50").Value = "  -2.66%  "

# Row 51
$ws.Range("E51").Value = "  -2.29%  "
